# Auto-generated edit script applying numeric corrections to market-price
# columns (H-N) across several worksheets, per upstream data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 269.75
$ws.Range("I6").Value = 304.7143
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 914.1428999999999
$ws.Range("L6").Value = 75
$ws.Range("M6").Value = -802.1428999999999
$ws.Range("N6").Value = -299

$ws.Range("H81").Value = 42999.11
$ws.Range("J81").Value = 42999.11
$ws.Range("L81").Value = 42999.11
$ws.Range("N81").Value = -44995.11

$ws.Range("H84").Value = 42999.11
$ws.Range("J84").Value = 42999.11
$ws.Range("L84").Value = 128997.33
$ws.Range("N84").Value = -138981.33

$ws.Range("H88").Value = 2354062.8
$ws.Range("I88").Value = 1350.7142
$ws.Range("K88").Value = 1350.7142
$ws.Range("M88").Value = -944.7141999999999

$ws.Range("H91").Value = 2354062.8
$ws.Range("I91").Value = 1350.7142
$ws.Range("K91").Value = 1350.7142
$ws.Range("M91").Value = 53.28580000000011

$ws.Range("H94").Value = 6000
$ws.Range("J94").Value = 6000
$ws.Range("L94").Value = 6000
$ws.Range("N94").Value = -6902

$ws.Range("H112").Value = 2855.4783
$ws.Range("J112").Value = 2932.6667
$ws.Range("L112").Value = 8798.000100000001
$ws.Range("N112").Value = -11014.0001

$ws.Range("H138").Value = 2783.7737
$ws.Range("J138").Value = 3202.279
$ws.Range("L138").Value = 9606.837
$ws.Range("N138").Value = -19886.837

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3495.725
$ws.Range("I32").Value = 3495.725
$ws.Range("K32").Value = 3495.725
$ws.Range("M32").Value = -3208.725

$ws.Range("H45").Value = 4072
$ws.Range("I45").Value = 3527.5
$ws.Range("K45").Value = 3527.5
$ws.Range("M45").Value = -3150.5

$ws.Range("H63").Value = 3088.8572
$ws.Range("I63").Value = 3267.4285
$ws.Range("K63").Value = 3267.4285
$ws.Range("M63").Value = -2581.4285

$ws.Range("H66").Value = 3088.8572
$ws.Range("I66").Value = 3267.4285
$ws.Range("K66").Value = 16337.1425
$ws.Range("M66").Value = -12905.1425

$ws.Range("H122").Value = 1444.625
$ws.Range("I122").Value = 1444.625
$ws.Range("K122").Value = 4333.875
$ws.Range("M122").Value = -1883.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = ""

$ws.Range("H31").Value = 1896.25
$ws.Range("I31").Value = 1430.75
$ws.Range("J31").Value = 2361.75
$ws.Range("K31").Value = 1430.75
$ws.Range("L31").Value = 2361.75
$ws.Range("M31").Value = -1135.75
$ws.Range("N31").Value = -2951.75

$ws.Range("H34").Value = 1896.25
$ws.Range("I34").Value = 1430.75
$ws.Range("J34").Value = 2361.75
$ws.Range("K34").Value = 1430.75
$ws.Range("L34").Value = 2361.75
$ws.Range("M34").Value = -1228.75
$ws.Range("N34").Value = -2765.75

$ws.Range("H70").Value = 79998.5
$ws.Range("J70").Value = 79998.5
$ws.Range("L70").Value = 79998.5
$ws.Range("N70").Value = -80628.5

$ws.Range("H73").Value = 79998.5
$ws.Range("J73").Value = 79998.5
$ws.Range("L73").Value = 79998.5
$ws.Range("N73").Value = -82182.5

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = ""
$ws.Range("N94").Value = ""

$ws.Range("H99").Value = 1607.421
$ws.Range("I99").Value = 1285.7142
$ws.Range("K99").Value = 1285.7142
$ws.Range("M99").Value = 212.2858000000001

$ws.Range("H107").Value = 1055.7693
$ws.Range("J107").Value = 611
$ws.Range("L107").Value = 611
$ws.Range("N107").Value = -4451

$ws.Range("H126").Value = 1607.421
$ws.Range("I126").Value = 1285.7142
$ws.Range("K126").Value = 3857.1426
$ws.Range("M126").Value = -1387.1426

$ws.Range("H132").Value = 13799.2
$ws.Range("I132").Value = 12249
$ws.Range("K132").Value = 36747
$ws.Range("M132").Value = -34217

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 4871.75
$ws.Range("J47").Value = 5830.6924
$ws.Range("L47").Value = 17492.0772
$ws.Range("N47").Value = -18354.0772

$ws.Range("H68").Value = 930.5
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""

$ws.Range("H71").Value = 930.5
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""

$ws.Range("H75").Value = 1950.8334
$ws.Range("I75").Value = 995
$ws.Range("J75").Value = 2428.75
$ws.Range("K75").Value = 2985
$ws.Range("L75").Value = 7286.25
$ws.Range("M75").Value = -1987
$ws.Range("N75").Value = -9282.25

$ws.Range("H78").Value = 1950.8334
$ws.Range("I78").Value = 995
$ws.Range("J78").Value = 2428.75
$ws.Range("K78").Value = 8955
$ws.Range("L78").Value = 21858.75
$ws.Range("M78").Value = -3963
$ws.Range("N78").Value = -31842.75

$ws.Range("H80").Value = 3972.5
$ws.Range("J80").Value = 4130
$ws.Range("L80").Value = 12390
$ws.Range("N80").Value = -14262

$ws.Range("H82").Value = 6013
$ws.Range("I82").Value = 6013
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 18039
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -17633
$ws.Range("N82").Value = ""

$ws.Range("H83").Value = 3972.5
$ws.Range("J83").Value = 4130
$ws.Range("L83").Value = 37170
$ws.Range("N83").Value = -46530

$ws.Range("H85").Value = 6013
$ws.Range("I85").Value = 6013
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 18039
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -16635
$ws.Range("N85").Value = ""

$ws.Range("H117").Value = 3767.8
$ws.Range("I117").Value = 2176.3333
$ws.Range("J117").Value = 4165.6665
$ws.Range("K117").Value = 6528.999899999999
$ws.Range("L117").Value = 12496.9995
$ws.Range("M117").Value = -3086.999899999999
$ws.Range("N117").Value = -19380.9995

$ws.Range("H131").Value = 536057.1
$ws.Range("I131").Value = 1011.7273
$ws.Range("K131").Value = 3035.1819
$ws.Range("M131").Value = 2004.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1533.4445
$ws.Range("I102").Value = 1466.0714
$ws.Range("K102").Value = 1466.0714
$ws.Range("M102").Value = 155.9286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3571
$ws.Range("I7").Value = 3332.8333
$ws.Range("K7").Value = 3332.8333
$ws.Range("M7").Value = -3220.8333

$ws.Range("H126").Value = 3571
$ws.Range("I126").Value = 3332.8333
$ws.Range("K126").Value = 9998.499899999999
$ws.Range("M126").Value = -7528.499899999999
